$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the (first) paragraph whose trimmed text exactly matches
# $text. Walking Paragraphs.Item(i) is more reliable here than Find.Execute,
# whose returned Range does not always expand back out to the full
# paragraph bounds.
# ---------------------------------------------------------------------------
function Find-ParaByText([string]$text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $cand = $d.Paragraphs.Item($i)
        if ($cand.Range.Text.Trim() -eq $text) {
            return $cand
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# Helper: insert a new list-paragraph immediately after $anchor, set its text
# and its bullet indent level (COM ListLevelNumber is 1-based; OOXML w:ilvl
# is 0-based, so ListLevelNumber = ilvl + 1). Returns the freshly created
# paragraph so callers can chain further inserts after it.
#
# NOTE: we deliberately re-look-up the new paragraph via its numeric Index
# (anchor's Index + 1) rather than $anchor.Next() -- when $anchor is the very
# last paragraph in the document body, Next() hands back a paragraph whose
# Range has no valid Start/End yet and silently drops any text assignment.
# Re-fetching by Paragraphs.Item(index) does not have that problem.
# ---------------------------------------------------------------------------
function Insert-ListItem($anchor, [string]$text, [int]$ilvl) {
    $idx = $anchor.Index
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($idx + 1)
    $newPara.Range.Text = $text
    $newPara.Range.ListFormat.ListLevelNumber = $ilvl + 1
    return $newPara
}

# ---------------------------------------------------------------------------
# 1) After "Site" (end of Timothy Couch's GUI sub-list) add the Prototype
#    sub-section with the new JMenuBar bullet.
# ---------------------------------------------------------------------------
$p = Find-ParaByText "Site"
$p = Insert-ListItem $p "Prototype" 1
$p = Insert-ListItem $p "Hierarchy of emails on left" 2
$p = Insert-ListItem $p "Main message panel" 2
$p = Insert-ListItem $p "Top JMenuBar thing" 2

# ---------------------------------------------------------------------------
# 2) After "User" (end of Will Hildreth's Account/User sub-list) add the
#    Sequence Diagrams / Close Email bullets.
# ---------------------------------------------------------------------------
$p = Find-ParaByText "User"
$p = Insert-ListItem $p "Sequence Diagrams" 1
$p = Insert-ListItem $p "Close Email" 2

# ---------------------------------------------------------------------------
# 3) After "Mailbox" (end of Daniel Johnson's GUI sub-list, the last
#    paragraph in the document) add another Prototype / Sequence Diagrams
#    block.
# ---------------------------------------------------------------------------
$p = Find-ParaByText "Mailbox"
$p = Insert-ListItem $p "Prototype" 1
$p = Insert-ListItem $p "Top row of mail buttons" 2
$p = Insert-ListItem $p "Sequence Diagrams" 1
$p = Insert-ListItem $p "Close Email" 2

Write-Host "Final paragraph count:" $d.Paragraphs.Count
